$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.493.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5085"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3855"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08458"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.60%  "

$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.495"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.818.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.65%  "

$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06708"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.07%  "

$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.513.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.275"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.026.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.400"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.092"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.07%  "

$ws.Range("E32").Value = "  -3.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.741"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.691"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07368"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2229"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.24%  "

$ws.Range("E37").Value = "  +0.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.223"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.786"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6320"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.192"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.61%  "

$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.759"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5904"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.194"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06983"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.47%  "
